# ============================================================================
# PlayerPerformance_4516.xlsx - additional scraping update
#  - adds a "Player Info" sheet (first tab)
#  - renames MATCH_CARD_LINK -> MATCH_CODE and stores the bare numeric match
#    code instead of the full scorecard URL on both "ODI Batting" and
#    "ODI Bowling"
#  - removes the stray empty INNING_NUMBER cells on "ODI Batting" for the
#    "did not bat" rows
#  - adds an "ODI Batting Extra" sheet (last tab) with additional batting
#    stats per match
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Helper: apply the bold / centred / top-aligned / thin-bordered header style
# that the original workbook uses on row 1 of every data sheet.
# ----------------------------------------------------------------------------
function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1         # xlContinuous (thin)
}

# ----------------------------------------------------------------------------
# 1. "Player Info" sheet - inserted before the current first sheet
# ----------------------------------------------------------------------------
$battingSheetRef = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetRef)
$playerInfo.Name = "Player Info"

# NOTE: the worksheet reference passed as the "Before" argument above gets
# reseated to the newly inserted sheet by this host's Add() binding, so the
# original "ODI Batting" sheet must be looked up again by name rather than
# reusing $battingSheetRef.
$battingSheet = $wb.Worksheets.Item("ODI Batting")

$playerInfo.Range("A1:D2").NumberFormat = "@"

$playerInfo.Cells.Item(1,1).Value = "ID"
$playerInfo.Cells.Item(1,2).Value = "NAME"
$playerInfo.Cells.Item(1,3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1,4).Value = "BOWL_STYLE"
Set-HeaderStyle $playerInfo.Range("A1:D1")

$playerInfo.Cells.Item(2,1).Value = "4516"
$playerInfo.Cells.Item(2,2).Value = "Mosaddeck Hossain"
$playerInfo.Cells.Item(2,3).Value = "Right Handed"
$playerInfo.Cells.Item(2,4).Value = "Right Arm Off Break"

# ----------------------------------------------------------------------------
# 2. "ODI Batting" - rename column header + collapse URL to bare match code
# ----------------------------------------------------------------------------
$battingSheet.Cells.Item(1,4).Value = "MATCH_CODE"

$battingLastRow = 44
$battingSheet.Range("D2:D" + $battingLastRow).NumberFormat = "@"
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $val = $cell.Value()
    if ($val -match "MatchCode=(\d+)") {
        $cell.Value = $matches[1]
    }
}

# Rows where INNING_NUMBER (column B) was a stray empty cell - remove them
# entirely so no <c> element remains for that coordinate.
$emptyInningRows = @(11, 13, 15, 20, 31, 41, 42, 44)
foreach ($r in $emptyInningRows) {
    $battingSheet.Cells.Item($r, 2).ClearContents()
}

# ----------------------------------------------------------------------------
# 3. "ODI Bowling" - rename column header + collapse URL to bare match code
# ----------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1,2).Value = "MATCH_CODE"

$bowlingLastRow = 42
$bowlingSheet.Range("B2:B" + $bowlingLastRow).NumberFormat = "@"
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $val = $cell.Value()
    if ($val -match "MatchCode=(\d+)") {
        $cell.Value = $matches[1]
    }
}

# ----------------------------------------------------------------------------
# 4. "ODI Batting Extra" - inserted after the current last sheet
# ----------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

$extra.Cells.Item(1,1).Value = "MATCH_CODE"
$extra.Cells.Item(1,2).Value = "BATTING_POSITION"
$extra.Cells.Item(1,3).Value = "NUM_4"
$extra.Cells.Item(1,4).Value = "NUM_6"
$extra.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Cells.Item(1,6).Value = "MAN_OF_MATCH"
Set-HeaderStyle $extra.Range("A1:F1")

# Text-formatted columns. MATCH_CODE / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL
# hold numeric look-alike strings in the source workbook, not real numbers -
# only BATTING_POSITION is a genuine number.
$extra.Range("A2:A21").NumberFormat = "@"
$extra.Range("C2:D21").NumberFormat = "@"
$extra.Range("E2:E21").NumberFormat = "@"
# Give column B (BATTING_POSITION, a real number) an explicit style too, so
# that the blank rows still keep an (empty) cell at that coordinate instead
# of disappearing from the sheet entirely.
$extra.Range("B2:B21").NumberFormat = "General"

$extraData = @(
    @("4199", 7,     "0", "0", "6.94%",  "NO"),
    @("4295", $null, $null, $null, $null, "NO"),
    @("4296", 7,     "2", "5", "24.41%", "YES"),
    @("4307", $null, $null, $null, $null, "NO"),
    @("4311", $null, $null, $null, $null, "NO"),
    @("4314", 7,     "4", "0", "9.29%",  "NO"),
    @("4325", $null, $null, $null, $null, "NO"),
    @("4335", 7,     "4", "0", "13.36%", "NO"),
    @("4345", 6,     "0", "0", "1.05%",  "NO"),
    @("4349", 7,     "1", "0", "7.24%",  "NO"),
    @("4356", $null, $null, $null, $null, "NO"),
    @("4357", $null, $null, $null, $null, "NO"),
    @("4464", $null, $null, $null, $null, "NO"),
    @("4465", 5,     "3", "1", "26.98%", "NO"),
    @("4477", 5,     "0", "0", "1.81%",  "NO"),
    @("4479", 5,     "0", "0", "2.07%",  "NO"),
    @("4481", 8,     $null, $null, $null, "NO"),
    @("4611", 8,     $null, $null, $null, "NO"),
    @("4616", 6,     "1", "1", "7.82%",  "NO"),
    @("4626", 6,     $null, $null, $null, "NO")
)

$r = 2
foreach ($row in $extraData) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $v = $row[$c]
        if ($null -eq $v) {
            $extra.Cells.Item($r, $c + 1).Value = ""
        } else {
            $extra.Cells.Item($r, $c + 1).Value = $v
        }
    }
    $r = $r + 1
}
